# Apply cryptos list update (commit: "Updated cryptos list on Sun Jul 30 11:18:23 UTC 2023 with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2 carries the workbook default (unstyled) cell format; we reuse it below so that
# forcing text-number-format on the Price column does not leave any data cell with a
# different style than before the edit.
$defaultStyle = $ws.Range("B2").Style

# The "Price" column holds numeric-looking text (e.g. "29.335.77", "1.000") that must
# stay literal text. Pre-format the whole data range as Text so assigning .Value keeps
# the exact string (matching zeros, multiple dots, etc.) instead of Excel coercing it
# into a Number.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '29.335.77'
$ws.Range("E2").Value = '  +0.16%  '
$ws.Range("D3").Value = '1.878.37'
$ws.Range("E3").Value = '  +0.37%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").Value = '0.7238'
$ws.Range("E5").Value = '  +2.06%  '
$ws.Range("D6").Value = '242.90'
$ws.Range("E6").Value = '  +0.53%  '
$ws.Range("E7").Value = '  +0.09%  '
$ws.Range("D8").Value = '0.08015'
$ws.Range("E8").Value = '  +2.67%  '
$ws.Range("D9").Value = '0.3159'
$ws.Range("E9").Value = '  +2.09%  '
$ws.Range("D10").Value = '25.04'
$ws.Range("E10").Value = '  +0.12%  '
$ws.Range("D11").Value = '0.08218'
$ws.Range("E11").Value = '  -2.24%  '
$ws.Range("D12").Value = '1.883.76'
$ws.Range("E12").Value = '  +0.51%  '
$ws.Range("D13").Value = '94.79'
$ws.Range("D14").Value = '5.231'
$ws.Range("E14").Value = '  -0.19%  '
$ws.Range("D15").Value = '0.7135'
$ws.Range("E15").Value = '  +0.47%  '
$ws.Range("D16").Value = '6.414'
$ws.Range("E16").Value = '  +5.71%  '
$ws.Range("D17").Value = '0.000008520'
$ws.Range("E17").Value = '  +3.92%  '
$ws.Range("D18").Value = '29.341.16'
$ws.Range("D19").Value = '243.92'
$ws.Range("E19").Value = '  +1.74%  '
$ws.Range("D20").Value = '13.29'
$ws.Range("E20").Value = '  +0.56%  '
$ws.Range("D21").Value = '1.000'
$ws.Range("E21").Value = '  +0.02%  '
$ws.Range("D22").Value = '7.781'
$ws.Range("E22").Value = '  +0.31%  '
$ws.Range("E23").Value = '  +0.05%  '
$ws.Range("D24").Value = '0.1608'
$ws.Range("E24").Value = '  +0.98%  '
$ws.Range("B25").Value = 'Cosmos'
$ws.Range("C25").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D25").Value = '9.052'
$ws.Range("E25").Value = '  +0.59%  '
$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D26").Value = '162.64'
$ws.Range("E26").Value = '  -0.40%  '
$ws.Range("D27").Value = '18.54'
$ws.Range("E27").Value = '  +0.51%  '
$ws.Range("D28").Value = '1.503'
$ws.Range("E28").Value = '  -0.11%  '
$ws.Range("D29").Value = '4.412'
$ws.Range("E29").Value = '  +0.55%  '
$ws.Range("D30").Value = '4.311'
$ws.Range("E30").Value = '  +0.53%  '
$ws.Range("D31").Value = '1.193'
$ws.Range("E31").Value = '  -7.88%  '
$ws.Range("D32").Value = '0.05376'
$ws.Range("E32").Value = '  +0.03%  '
$ws.Range("D33").Value = '1.938'
$ws.Range("E33").Value = '  -0.04%  '
$ws.Range("E34").Value = '  +1.77%  '
$ws.Range("D35").Value = '1.179'
$ws.Range("E35").Value = '  +0.25%  '
$ws.Range("D36").Value = '2.707'
$ws.Range("E36").Value = '  +0.45%  '
$ws.Range("B37").Value = 'Maker'
$ws.Range("C37").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D37").Value = '1.283.46'
$ws.Range("E37").Value = '  +4.52%  '
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").Value = '0.01877'
$ws.Range("E38").Value = '  +0.41%  '
$ws.Range("D39").Value = '2.754'
$ws.Range("E39").Value = '  +1.01%  '
$ws.Range("D40").Value = '6.450'
$ws.Range("E40").Value = '  -1.47%  '
$ws.Range("D41").Value = '113.30'
$ws.Range("E41").Value = '  +4.40%  '
$ws.Range("D42").Value = '0.9112'
$ws.Range("E42").Value = '  +2.89%  '
$ws.Range("D43").Value = '74.45'
$ws.Range("E43").Value = '  +2.83%  '
$ws.Range("E44").Value = '  +8.62%  '
$ws.Range("E45").Value = '  +0.11%  '
$ws.Range("D46").Value = '2.029.10'
$ws.Range("E46").Value = '  +0.40%  '
$ws.Range("D47").Value = '0.5225'
$ws.Range("E47").Value = '  +0.61%  '
$ws.Range("D48").Value = '1.798'
$ws.Range("E48").Value = '  +0.50%  '
$ws.Range("D49").Value = '9.501'
$ws.Range("E49").Value = '  +0.97%  '
$ws.Range("D50").Value = '0.4353'
$ws.Range("E50").Value = '  +1.01%  '
$ws.Range("D51").Value = '7.119'
$ws.Range("E51").Value = '  +0.68%  '

# Restore the original (default) style on the whole range so no stray "s" attribute
# or extra referenced style is introduced on the price/volume cells.
$ws.Range("D2:E51").Style = $defaultStyle

